$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings like "43.806.57" or "239.50" that
# must stay literal text (dotted thousands separators, trailing zeros).
# Force text format before assigning so Excel does not coerce them to
# numbers and strip the formatting.

$ws.Range("D2").Value = "43.806.57"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.346.74"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.50"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.667"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.40"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.95"
$ws.Range("E11").Value = "  +6.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.59"
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.24"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.14"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.901"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "2.353.34"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "43.728.05"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "252.79"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.35"
$ws.Range("E27").Value = "  -3.65%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.76"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.21"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.134"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0740"
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("E34").Value = "  -4.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.35"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.40"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0271"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.31"
$ws.Range("E40").Value = "  +12.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.59"
$ws.Range("E41").Value = "  +14.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.49"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.13"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.105"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.200"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.84"
$ws.Range("E50").Value = "  -3.00%  "
$ws.Range("E51").Value = "  +2.62%  "
